# Reproduce the worksheet-view / column-width part of the
# "bug fix on time series writing" commit on the "Electricity" sheet:
#  - the view no longer sits scrolled down to A40 with A70:XFD70
#    selected; it is back near the top with a single cell, I9, selected
#  - column B (country / node names, e.g. "Bosnia and Hertz",
#    "Northern Ireland") is widened to fit its longest entry

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Electricity")
[void]$ws.Activate()

# Best-fit column B to its contents (Format > Column Width > AutoFit
# Selection), then pin it to the width recorded in the saved workbook.
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(2).ColumnWidth = 13

# Move the selection up to I9 - this also clears the old topLeftCell="A40"
# scroll position, since I9 is already within the default viewport.
$ws.Range("I9").Select() | Out-Null
